# Add data for 2022-04-24 (update report "through" date from 04-15 to 04-16)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-16"

# Update the April row label text
$ws.Range("A5").Value = "April (through 04-16)"

# February 2022 (column I, row 3)
$ws.Range("I3").Value = 141

# April row (row 5): columns C..I updated, B unchanged
$ws.Range("C5").Value = 12
$ws.Range("D5").Value = 31
$ws.Range("E5").Value = 28
$ws.Range("F5").Value = 27
$ws.Range("G5").Value = 34
$ws.Range("H5").Value = 55
$ws.Range("I5").Value = 68

# Total row (row 6): columns C..I updated, B unchanged
$ws.Range("C6").Value = 140
$ws.Range("D6").Value = 220
$ws.Range("E6").Value = 225
$ws.Range("F6").Value = 137
$ws.Range("G6").Value = 232
$ws.Range("H6").Value = 478
$ws.Range("I6").Value = 503
